$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.514.69"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.36"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.41"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4815"
$ws.Range("E7").Value = "  +2.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2891"
$ws.Range("E8").Value = "  +1.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06729"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.79"
$ws.Range("E10").Value = "  +2.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.25"
$ws.Range("E11").Value = "  +5.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.919.80"
$ws.Range("E12").Value = "  +1.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07561"
$ws.Range("E13").Value = "  -2.14%  "

$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6716"
$ws.Range("E15").Value = "  +2.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "287.66"
$ws.Range("E16").Value = "  -2.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.519.56"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007607"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.167.77"
$ws.Range("E21").Value = "  +1.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.495"
$ws.Range("E22").Value = "  +5.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.416"
$ws.Range("E24").Value = "  +2.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.475"
$ws.Range("E25").Value = "  +1.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.36"
$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.35"
$ws.Range("E27").Value = "  -5.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.141"
$ws.Range("E28").Value = "  +2.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1063"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("E30").Value = "  +2.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.165"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.035"
$ws.Range("E32").Value = "  +1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04992"
$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7296"
$ws.Range("E34").Value = "  -1.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02053"
$ws.Range("E36").Value = "  -1.78%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.735"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.670"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.83"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("E41").Value = "  -2.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4441"
$ws.Range("E42").Value = "  +4.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8653"
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.888"
$ws.Range("E44").Value = "  +1.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.19"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("E47").Value = "  +2.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.88"
$ws.Range("E48").Value = "  -4.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.342"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1243"
$ws.Range("E50").Value = "  +2.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.84"
$ws.Range("E51").Value = "  +0.10%  "
